$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "SN1/SN2/SN3" fatty-acid chain labels were renamed to "FA1/FA2/FA3"
# (figure TG [M+Na]+ / [M+H]+ weighting labels) across column A, rows 2-13.
$replacements = @{
    "A2"  = "FA1_[FA-H2O+H]+"
    "A3"  = "FA2_[FA-H2O+H]+"
    "A4"  = "FA3_[FA-H2O+H]+"
    "A5"  = "[MG(FA1)-H2O+H]+"
    "A6"  = "[MG(FA2)-H2O+H]+"
    "A7"  = "[MG(FA3)-H2O+H]+"
    "A8"  = "[M-(FA1)+Na]+"
    "A9"  = "[M-(FA2)+Na]+"
    "A10" = "[M-(FA3)+Na]+"
    "A11" = "[M-(FA1-H+Na)+H]+"
    "A12" = "[M-(FA2-H+Na)+H]+"
    "A13" = "[M-(FA3-H+Na)+H]+"
}

foreach ($addr in $replacements.Keys) {
    $ws.Range($addr).Value = $replacements[$addr]
}

# Scroll the view down and move the selection to A13, matching the saved
# window state (top visible row 7, active cell A13).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A13").Select()
